# Insert two new weekly price rows for "Sandia" (Vega Modelo de Temuco) right
# after the existing row 655, pushing all subsequent rows down by two and
# growing the used range from A1:R756 to A1:R758.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 656-657 (existing rows 656..756 shift to 658..758).
$ws.Range("A656:R657").EntireRow.Insert()

# New row 656: Primera, Peru, $/kilo (volumen en unidades)
$ws.Cells.Item(656, 1).Value = 10
$ws.Cells.Item(656, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(656, 3).Value = "La Araucanía"
$ws.Cells.Item(656, 4).Value = 45212
$ws.Cells.Item(656, 5).Value = 9
$ws.Cells.Item(656, 6).Value = 100112028
$ws.Cells.Item(656, 7).Value = "Sandia"
$ws.Cells.Item(656, 8).Value = "Sin especificar"
$ws.Cells.Item(656, 9).Value = "Primera"
$ws.Cells.Item(656, 10).Value = 300
$ws.Cells.Item(656, 11).Value = 900
$ws.Cells.Item(656, 12).Value = 900
$ws.Cells.Item(656, 13).Value = 900
$ws.Cells.Item(656, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(656, 15).Value = "Perú"
$ws.Cells.Item(656, 16).Value = 900
$ws.Cells.Item(656, 17).Value = 1
$ws.Cells.Item(656, 18).Value = "Hortaliza"

# New row 657: Primera, Brasil, $/unidad
$ws.Cells.Item(657, 1).Value = 10
$ws.Cells.Item(657, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(657, 3).Value = "La Araucanía"
$ws.Cells.Item(657, 4).Value = 45212
$ws.Cells.Item(657, 5).Value = 9
$ws.Cells.Item(657, 6).Value = 100112028
$ws.Cells.Item(657, 7).Value = "Sandia"
$ws.Cells.Item(657, 8).Value = "Sin especificar"
$ws.Cells.Item(657, 9).Value = "Primera"
$ws.Cells.Item(657, 10).Value = 220
$ws.Cells.Item(657, 11).Value = 3200
$ws.Cells.Item(657, 12).Value = 3200
$ws.Cells.Item(657, 13).Value = 3200
$ws.Cells.Item(657, 14).Value = "$/unidad"
$ws.Cells.Item(657, 15).Value = "Brasil"
$ws.Cells.Item(657, 16).Value = 3200
$ws.Cells.Item(657, 17).Value = 1
$ws.Cells.Item(657, 18).Value = "Hortaliza"
